# Updated cryptos list (Sun Oct 13 08:43:12 UTC 2024) via GitHub Actions.
# Source data is plain text (coin name / link / price / 1h volume %),
# so every write below stays a text value. For cells whose new value
# looks like a plain number (e.g. '574.34'), a leading apostrophe is
# used (classic Excel "force text" quote-prefix) so Excel keeps storing
# it as text instead of silently converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.894.83'
# Row 3
$ws.Range('D3').Value = '2.464.82'
$ws.Range('E3').Value = '  +0.76%  '
# Row 5
$ws.Range('D5').Value = '''574.34'
$ws.Range('E5').Value = '  -0.23%  '
# Row 6
$ws.Range('D6').Value = '''146.86'
$ws.Range('E6').Value = '  +0.76%  '
# Row 7
$ws.Range('E7').Value = '  -0.04%  '
# Row 8
$ws.Range('E8').Value = '  -0.39%  '
# Row 9
$ws.Range('D9').Value = '2.465.79'
# Row 10
$ws.Range('E10').Value = '  +0.43%  '
# Row 11
$ws.Range('D11').Value = '''0.162'
$ws.Range('E11').Value = '  -0.42%  '
# Row 12
$ws.Range('E12').Value = '  +0.74%  '
# Row 13
$ws.Range('E13').Value = '  +1.22%  '
# Row 14
$ws.Range('D14').Value = '''29.06'
$ws.Range('E14').Value = '  +3.04%  '
# Row 15
$ws.Range('E15').Value = '  -0.37%  '
# Row 16
$ws.Range('D16').Value = '2.912.02'
$ws.Range('E16').Value = '  +0.76%  '
# Row 17
$ws.Range('D17').Value = '62.755.16'
# Row 18
$ws.Range('D18').Value = '2.467.14'
$ws.Range('E18').Value = '  +0.97%  '
# Row 19
$ws.Range('D19').Value = '''7.94'
$ws.Range('E19').Value = '  +0.39%  '
# Row 20
$ws.Range('D20').Value = '''10.99'
$ws.Range('E20').Value = '  -0.06%  '
# Row 21
$ws.Range('D21').Value = '''327.20'
$ws.Range('E21').Value = '  -0.77%  '
# Row 22
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = '''4.13'
$ws.Range('E22').Value = '  -0.04%  '
# Row 23
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value = '''2.22'
$ws.Range('E23').Value = '  +8.69%  '
# Row 24
$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  -0.14%  '
# Row 25
$ws.Range('D25').Value = '''10.02'
$ws.Range('E25').Value = '  +17.37%  '
# Row 26
$ws.Range('D26').Value = '''65.53'
$ws.Range('E26').Value = '  -0.79%  '
# Row 27
$ws.Range('D27').Value = '''647.23'
$ws.Range('E27').Value = '  -0.54%  '
# Row 28
$ws.Range('D28').Value = '0.0₃0990'
$ws.Range('E28').Value = '  -0.05%  '
# Row 29
$ws.Range('E29').Value = '  +0.98%  '
# Row 30
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -15.08%  '
# Row 31
$ws.Range('D31').Value = '''1.43'
$ws.Range('E31').Value = '  -1.88%  '
# Row 32
$ws.Range('D32').Value = '''7.98'
$ws.Range('E32').Value = '  -2.84%  '
# Row 33
$ws.Range('E33').Value = '  -1.14%  '
# Row 34
$ws.Range('E34').Value = '  -3.27%  '
# Row 35
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  +0.00%  '
# Row 36
$ws.Range('D36').Value = '''1.54'
$ws.Range('E36').Value = '  +3.09%  '
# Row 37
$ws.Range('D37').Value = '''4.75'
$ws.Range('E37').Value = '  -0.37%  '
# Row 38
$ws.Range('D38').Value = '''2.85'
$ws.Range('E38').Value = '  +4.34%  '
# Row 39
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '''0.369'
$ws.Range('E39').Value = '  -1.46%  '
# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '''5.41'
$ws.Range('E40').Value = '  -1.96%  '
# Row 41
$ws.Range('D41').Value = '''151.46'
$ws.Range('E41').Value = '  -0.77%  '
# Row 42
$ws.Range('E42').Value = '  -0.33%  '
# Row 43
$ws.Range('E43').Value = '  -1.28%  '
# Row 44
$ws.Range('D44').Value = '0.0₆0308'
$ws.Range('E44').Value = '  -48.92%  '
# Row 45
$ws.Range('E45').Value = '  +0.00%  '
# Row 46
$ws.Range('D46').Value = '''152.75'
$ws.Range('E46').Value = '  +4.93%  '
# Row 47
$ws.Range('D47').Value = '''15.24'
$ws.Range('E47').Value = '  +2.06%  '
# Row 48
$ws.Range('E48').Value = '  -1.61%  '
# Row 49
$ws.Range('D49').Value = '''20.53'
$ws.Range('E49').Value = '  -1.22%  '
# Row 50
$ws.Range('E50').Value = '  +0.42%  '
# Row 51
$ws.Range('D51').Value = '''0.0511'
$ws.Range('E51').Value = '  -1.24%  '
